$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (column F) values to reflect repulled data / recalculated mean
$ws.Range("F2").Value = -1
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -2
$ws.Range("F13").Value = 2
$ws.Range("F15").Value = -2
